$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.880.14'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.462.37'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '2.462.63'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('E10').Value = '  +1.03%  '
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.356'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '2.908.32'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '62.772.51'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '2.463.52'
$ws.Range('E18').Value = '  +1.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '327.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.24'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.96%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +20.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.77'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '649.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.03%  '
$ws.Range('D28').Value = '0.0₃0983'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '2.584.59'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('E30').Value = '  -13.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.45'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.85'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('E34').Value = '  -3.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.369'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '151.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('D44').Value = '0.0₆0317'
$ws.Range('E44').Value = '  -61.84%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '153.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.00%  '
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.45'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.608'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0512'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.38%  '
